$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.042103160828266
$ws.Range("D2").Value = 1.047449668333418
$ws.Range("E2").Value = 1.045762625991648
$ws.Range("F2").Value = 1.055241864462713
$ws.Range("I2").Value = 1.036080199957325
$ws.Range("J2").Value = 1.047180785994238
$ws.Range("K2").Value = 1.050212273683783
$ws.Range("L2").Value = 1.048529955483006
$ws.Range("M2").Value = 1.057982862116065
$ws.Range("N2").Value = 1.019437771148988

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.043592275506503
$ws.Range("D3").Value = 1.048622793244932
$ws.Range("E3").Value = 1.047201647717809
$ws.Range("F3").Value = 1.056633465450365
$ws.Range("I3").Value = 1.036393882336288
$ws.Range("J3").Value = 1.048313519075403
$ws.Range("K3").Value = 1.051196308484739
$ws.Range("L3").Value = 1.049778851263619
$ws.Range("M3").Value = 1.059186390690868
$ws.Range("N3").Value = 1.019829062078323

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.044554474257271
$ws.Range("D4").Value = 1.049380465269363
$ws.Range("E4").Value = 1.048131808629375
$ws.Range("F4").Value = 1.057532866797813
$ws.Range("I4").Value = 1.036594861667678
$ws.Range("J4").Value = 1.049044731523641
$ws.Range("K4").Value = 1.05183105651729
$ws.Range("L4").Value = 1.05058547547313
$ws.Range("M4").Value = 1.059963569574531
$ws.Range("N4").Value = 1.020081282631951

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.044958664855833
$ws.Range("D5").Value = 1.049698655890445
$ws.Range("E5").Value = 1.048522619470004
$ws.Range("F5").Value = 1.057910727795737
$ws.Range("I5").Value = 1.036678877663116
$ws.Range("J5").Value = 1.049351721638199
$ws.Range("K5").Value = 1.05209743310303
$ws.Range("L5").Value = 1.050924228172968
$ws.Range("M5").Value = 1.060289921542084
$ws.Range("N5").Value = 1.020187085434746

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.045026511815536
$ws.Range("D6").Value = 1.049752062008785
$ws.Range("E6").Value = 1.048588225097062
$ws.Range("F6").Value = 1.057974157985858
$ws.Range("I6").Value = 1.036692956471018
$ws.Range("J6").Value = 1.049403242647601
$ws.Range("K6").Value = 1.052142131386228
$ws.Range("L6").Value = 1.050981085760186
$ws.Range("M6").Value = 1.060344695693761
$ws.Range("N6").Value = 1.020204836710758

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.044559876311293
$ws.Range("D7").Value = 1.049384718257186
$ws.Range("E7").Value = 1.048137031554121
$ws.Range("F7").Value = 1.057537916759559
$ws.Range("I7").Value = 1.036595986160684
$ws.Range("J7").Value = 1.049048835151281
$ws.Range("K7").Value = 1.051834617701326
$ws.Range("L7").Value = 1.050590003279752
$ws.Range("M7").Value = 1.059967931768631
$ws.Range("N7").Value = 1.020082697278641

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.042606699731954
$ws.Range("D8").Value = 1.047846427962423
$ws.Range("E8").Value = 1.046249157374935
$ws.Range("F8").Value = 1.055712385034782
$ws.Range("I8").Value = 1.036186624305382
$ws.Range("J8").Value = 1.04756396227764
$ws.Range("K8").Value = 1.050545247263836
$ws.Range("L8").Value = 1.048952338657612
$ws.Range("M8").Value = 1.058389931540743
$ws.Range("N8").Value = 1.019570211799098

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.039154221699902
$ws.Range("D9").Value = 1.045124684369643
$ws.Range("E9").Value = 1.042914671487587
$ws.Range("F9").Value = 1.052487202210424
$ws.Range("I9").Value = 1.03544992754739
$ws.Range("J9").Value = 1.044933847610834
$ws.Range("K9").Value = 1.04825778076856
$ws.Range("L9").Value = 1.046054856488974
$ws.Range("M9").Value = 1.055596918265211
$ws.Range("N9").Value = 1.018659634593932

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.036844890421676
$ws.Range("D10").Value = 1.043302432616796
$ws.Range("E10").Value = 1.040686026682321
$ws.Range("F10").Value = 1.050331087313904
$ws.Range("I10").Value = 1.03494836888239
$ws.Range("J10").Value = 1.04317096985355
$ws.Range("K10").Value = 1.046722141756622
$ws.Range("L10").Value = 1.044114965761424
$ws.Range("M10").Value = 1.05372625206327
$ws.Range("N10").Value = 1.01804742372127

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.035842999226458
$ws.Range("D11").Value = 1.042511468828855
$ws.Range("E11").Value = 1.039719567449367
$ws.Range("F11").Value = 1.049395961136644
$ws.Range("I11").Value = 1.034728690771112
$ws.Range("J11").Value = 1.042405302728157
$ws.Range("K11").Value = 1.046054599955737
$ws.Range("L11").Value = 1.04327293649937
$ws.Range("M11").Value = 1.052914103847109
$ws.Range("N11").Value = 1.017781081411866

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.035470552414972
$ws.Range("D12").Value = 1.042217375646151
$ws.Range("E12").Value = 1.039360357052154
$ws.Range("F12").Value = 1.049048378892083
$ws.Range("I12").Value = 1.034646714789643
$ws.Range("J12").Value = 1.042120543269467
$ws.Range("K12").Value = 1.045806248947295
$ws.Range("L12").Value = 1.042959855161059
$ws.Range("M12").Value = 1.052612108103886
$ws.Range("N12").Value = 1.017681959919823

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.035550457233074
$ws.Range("D13").Value = 1.042280473034502
$ws.Range("E13").Value = 1.039437419183496
$ws.Range("F13").Value = 1.049122947203274
$ws.Range("I13").Value = 1.034664316044931
$ws.Range("J13").Value = 1.042181641383725
$ws.Range("K13").Value = 1.045859539104434
$ws.Range("L13").Value = 1.043027026503157
$ws.Range("M13").Value = 1.052676902153402
$ws.Range("N13").Value = 1.017703230453574

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.035812218825354
$ws.Range("D14").Value = 1.042487165027572
$ws.Range("E14").Value = 1.039689879629406
$ws.Range("F14").Value = 1.049367234710794
$ws.Range("I14").Value = 1.034721922329503
$ws.Range("J14").Value = 1.042381771735466
$ws.Range("K14").Value = 1.046034079300334
$ws.Range("L14").Value = 1.04324706354885
$ws.Range("M14").Value = 1.052889147508905
$ws.Range("N14").Value = 1.017772891888928

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.035973458782443
$ws.Range("D15").Value = 1.042614475644799
$ws.Range("E15").Value = 1.039845398870059
$ws.Range("F15").Value = 1.049517716982575
$ws.Range("I15").Value = 1.034757365341974
$ws.Range("J15").Value = 1.042505031198587
$ws.Range("K15").Value = 1.046141566597459
$ws.Range("L15").Value = 1.043382593788974
$ws.Range("M15").Value = 1.053019875250603
$ws.Range("N15").Value = 1.017815787343501

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.036911341092951
$ws.Range("D16").Value = 1.043354885379924
$ws.Range("E16").Value = 1.040750136361759
$ws.Range("F16").Value = 1.050393116038188
$ws.Range("I16").Value = 1.03496289534051
$ws.Range("J16").Value = 1.043221735016225
$ws.Range("K16").Value = 1.046766389048244
$ws.Range("L16").Value = 1.044170804803506
$ws.Range("M16").Value = 1.053780106069629
$ws.Range("N16").Value = 1.018065073443206

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.037499125346336
$ws.Range("D17").Value = 1.043818807424783
$ws.Range("E17").Value = 1.041317262870668
$ws.Range("F17").Value = 1.050941820239988
$ws.Range("I17").Value = 1.035091148059692
$ws.Range("J17").Value = 1.043670676070352
$ws.Range("K17").Value = 1.047157623606307
$ws.Range("L17").Value = 1.044664676558906
$ws.Range("M17").Value = 1.054256401880119
$ws.Range("N17").Value = 1.018221107592932

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.037841783821397
$ws.Range("D18").Value = 1.044089220777217
$ws.Range("E18").Value = 1.041647919636204
$ws.Range("F18").Value = 1.051261724513025
$ws.Range("I18").Value = 1.035165714644665
$ws.Range("J18").Value = 1.043932311487111
$ws.Range("K18").Value = 1.04738557362926
$ws.Range("L18").Value = 1.044952546888187
$ws.Range("M18").Value = 1.054534011389523
$ws.Range("N18").Value = 1.018311999119894

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.03795859020623
$ws.Range("D19").Value = 1.044181393625004
$ws.Range("E19").Value = 1.041760641787564
$ws.Range("F19").Value = 1.051370779124297
$ws.Range("I19").Value = 1.035191099101001
$ws.Range("J19").Value = 1.044021484556454
$ws.Range("K19").Value = 1.047463256381111
$ws.Range("L19").Value = 1.045050670080243
$ws.Range("M19").Value = 1.054628634353532
$ws.Range("N19").Value = 1.018342970398742

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.037436080976938
$ws.Range("D20").Value = 1.043769052110258
$ws.Range("E20").Value = 1.041256429949534
$ws.Range("F20").Value = 1.050882964585195
$ws.Range("I20").Value = 1.035077412706009
$ws.Range("J20").Value = 1.043622532192911
$ws.Range("K20").Value = 1.047115673789635
$ws.Range("L20").Value = 1.044611709191802
$ws.Range("M20").Value = 1.054205321179878
$ws.Range("N20").Value = 1.018204379090333

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.035735144909753
$ws.Range("D21").Value = 1.042426307552294
$ws.Range("E21").Value = 1.039615542578668
$ws.Range("F21").Value = 1.049295304667216
$ws.Range("I21").Value = 1.034704969159746
$ws.Range("J21").Value = 1.042322848232303
$ws.Range("K21").Value = 1.045982692547255
$ws.Range("L21").Value = 1.043182276853547
$ws.Range("M21").Value = 1.05282665564225
$ws.Range("N21").Value = 1.017752383597572

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.034663959754354
$ws.Range("D22").Value = 1.041580366649536
$ws.Range("E22").Value = 1.038582549276304
$ws.Range("F22").Value = 1.048295718060741
$ws.Range("I22").Value = 1.034468612488365
$ws.Range("J22").Value = 1.041503618861177
$ws.Range("K22").Value = 1.04526804679699
$ws.Range("L22").Value = 1.042281714696156
$ws.Range("M22").Value = 1.05195793441256
$ws.Range("N22").Value = 1.017467095152714

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.035231983152333
$ws.Range("D23").Value = 1.042028979487574
$ws.Range("E23").Value = 1.039130284652731
$ws.Range("F23").Value = 1.048825749352926
$ws.Range("I23").Value = 1.034594117627875
$ws.Range("J23").Value = 1.041938105955215
$ws.Range("K23").Value = 1.045647113511258
$ws.Range("L23").Value = 1.042759294766796
$ws.Range("M23").Value = 1.052418642243758
$ws.Range("N23").Value = 1.017618437036211

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.03746456858913
$ws.Range("D24").Value = 1.043791534965653
$ws.Range("E24").Value = 1.041283918156398
$ws.Range("F24").Value = 1.050909559372386
$ws.Range("I24").Value = 1.035083619866058
$ws.Range("J24").Value = 1.043644287032107
$ws.Range("K24").Value = 1.047134629880257
$ws.Range("L24").Value = 1.044635643473104
$ws.Range("M24").Value = 1.054228402985674
$ws.Range("N24").Value = 1.018211938353191

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.040048089332273
$ws.Range("D25").Value = 1.045829664740421
$ws.Range("E25").Value = 1.043777682894154
$ws.Range("F25").Value = 1.053322018263654
$ws.Range("I25").Value = 1.035642211273466
$ws.Range("J25").Value = 1.045615439139742
$ws.Range("K25").Value = 1.048851003688589
$ws.Range("L25").Value = 1.046805349858798
$ws.Range("M25").Value = 1.056320480205388
$ws.Range("N25").Value = 1.018895941991413

Write-Output "vm_pu values updated for 380 kV case (rows 2-25)"